# Update "Last Updated" timestamp on the Metadata sheet
$wbMeta = $excel.ActiveWorkbook
$wsMeta = $wbMeta.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value2 = "05 Nov 2025, 03:03 PM"

# Insert a new top row of stock data on the "Stock List" sheet, pushing
# the rest of the table down by one row (the former last row falls off
# the bottom of the table, matching the unchanged A1:H76 dimension).
$ws = $wbMeta.Worksheets.Item("Stock List")
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(77).Delete()

$ws.Range("A2").Value2 = "📋"
$ws.Range("B2").Value2 = "CAPTRU-RE1"
$ws.Range("C2").Value2 = "CAPTRU-RE1"
$ws.Range("D2").Value2 = 5.67
$ws.Range("E2").Value2 = -11.9565
$ws.Range("F2").Value2 = "N/A"
$ws.Range("G2").Value2 = "N/A"
$ws.Range("H2").Value2 = 0
